# Weekly update: insert a new data row for "Macroferia Regional de Talca -
# Espárragos" before current row 78, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78; this pushes rows 78:108 down to 79:109
# and keeps the D-column date style (s="2") that Excel carries over from
# the row being pushed down.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 45202
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = 300000000
$ws.Range("G78").Value = "Espárragos"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 1100
$ws.Range("L78").Value = 1200
$ws.Range("M78").Value = 1140
$ws.Range("N78").Value = "$/kilo"
$ws.Range("O78").Value = "Provincia de Linares"
$ws.Range("P78").Value = 1140
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"
